$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) columns with 2020-08-31 data
$ws.Cells.Item(2, 3).Value = 38672
$ws.Cells.Item(2, 4).Value = 55919422
$ws.Cells.Item(3, 3).Value = 92751
$ws.Cells.Item(3, 4).Value = 135948203
$ws.Cells.Item(4, 3).Value = 31708
$ws.Cells.Item(4, 4).Value = 46954945
$ws.Cells.Item(5, 3).Value = 8899
$ws.Cells.Item(5, 4).Value = 13225897
$ws.Cells.Item(6, 3).Value = 2065
$ws.Cells.Item(6, 4).Value = 3069471
$ws.Cells.Item(7, 3).Value = 168
$ws.Cells.Item(7, 4).Value = 247093
$ws.Cells.Item(12, 3).Value = 42128
$ws.Cells.Item(12, 4).Value = 57124908
$ws.Cells.Item(13, 3).Value = 9867
$ws.Cells.Item(13, 4).Value = 14271278
$ws.Cells.Item(14, 3).Value = 26389
$ws.Cells.Item(14, 4).Value = 38688507
$ws.Cells.Item(15, 3).Value = 8432
$ws.Cells.Item(15, 4).Value = 12513478
$ws.Cells.Item(16, 3).Value = 2199
$ws.Cells.Item(16, 4).Value = 3268153
$ws.Cells.Item(17, 3).Value = 430
$ws.Cells.Item(17, 4).Value = 634123
$ws.Cells.Item(20, 3).Value = 10384
$ws.Cells.Item(20, 4).Value = 13737829
$ws.Cells.Item(21, 3).Value = 13640
$ws.Cells.Item(21, 4).Value = 19685262
$ws.Cells.Item(22, 3).Value = 32158
$ws.Cells.Item(22, 4).Value = 47183212
$ws.Cells.Item(23, 3).Value = 10375
$ws.Cells.Item(23, 4).Value = 15421110
$ws.Cells.Item(24, 3).Value = 2689
$ws.Cells.Item(24, 4).Value = 3998271
$ws.Cells.Item(25, 3).Value = 529
$ws.Cells.Item(25, 4).Value = 787592
$ws.Cells.Item(27, 3).Value = 11896
$ws.Cells.Item(27, 4).Value = 15879250
$ws.Cells.Item(28, 3).Value = 7847
$ws.Cells.Item(28, 4).Value = 11355707
$ws.Cells.Item(29, 3).Value = 22934
$ws.Cells.Item(29, 4).Value = 33664023
$ws.Cells.Item(30, 3).Value = 7921
$ws.Cells.Item(30, 4).Value = 11781491
$ws.Cells.Item(31, 3).Value = 2003
$ws.Cells.Item(31, 4).Value = 2988751
$ws.Cells.Item(32, 3).Value = 377
$ws.Cells.Item(32, 4).Value = 562915
$ws.Cells.Item(34, 3).Value = 8455
$ws.Cells.Item(34, 4).Value = 11165655
$ws.Cells.Item(35, 3).Value = 3341
$ws.Cells.Item(35, 4).Value = 4824946
$ws.Cells.Item(36, 3).Value = 8002
$ws.Cells.Item(36, 4).Value = 11685621
$ws.Cells.Item(37, 3).Value = 3232
$ws.Cells.Item(37, 4).Value = 4790961
$ws.Cells.Item(38, 3).Value = 837
$ws.Cells.Item(38, 4).Value = 1246723
$ws.Cells.Item(41, 3).Value = 2526
$ws.Cells.Item(41, 4).Value = 3412272
$ws.Cells.Item(42, 3).Value = 17661
$ws.Cells.Item(42, 4).Value = 25535774
$ws.Cells.Item(43, 3).Value = 52089
$ws.Cells.Item(43, 4).Value = 76349533
$ws.Cells.Item(44, 3).Value = 19286
$ws.Cells.Item(44, 4).Value = 28641344
$ws.Cells.Item(45, 3).Value = 5721
$ws.Cells.Item(45, 4).Value = 8516760
$ws.Cells.Item(46, 3).Value = 1247
$ws.Cells.Item(46, 4).Value = 1861045
$ws.Cells.Item(50, 3).Value = 17084
$ws.Cells.Item(50, 4).Value = 22700975
$ws.Cells.Item(51, 3).Value = 2117
$ws.Cells.Item(51, 4).Value = 3071378
$ws.Cells.Item(52, 3).Value = 7167
$ws.Cells.Item(52, 4).Value = 10532201
$ws.Cells.Item(53, 3).Value = 2415
$ws.Cells.Item(53, 4).Value = 3606964
$ws.Cells.Item(57, 3).Value = 7290
$ws.Cells.Item(57, 4).Value = 10021130
$ws.Cells.Item(58, 3).Value = 1135
$ws.Cells.Item(58, 4).Value = 1893144
$ws.Cells.Item(59, 3).Value = 2785
$ws.Cells.Item(59, 4).Value = 4626488
$ws.Cells.Item(60, 3).Value = 1095
$ws.Cells.Item(60, 4).Value = 1820338
$ws.Cells.Item(61, 3).Value = 377
$ws.Cells.Item(61, 4).Value = 630883
$ws.Cells.Item(62, 3).Value = 123
$ws.Cells.Item(62, 4).Value = 209600
$ws.Cells.Item(64, 3).Value = 1660
$ws.Cells.Item(64, 4).Value = 2566499
$ws.Cells.Item(65, 3).Value = 15728
$ws.Cells.Item(65, 4).Value = 22713768
$ws.Cells.Item(66, 3).Value = 45518
$ws.Cells.Item(66, 4).Value = 66594018
$ws.Cells.Item(67, 3).Value = 15941
$ws.Cells.Item(67, 4).Value = 23684645
$ws.Cells.Item(68, 3).Value = 4638
$ws.Cells.Item(68, 4).Value = 6908051
$ws.Cells.Item(69, 3).Value = 959
$ws.Cells.Item(69, 4).Value = 1426668
$ws.Cells.Item(73, 3).Value = 15362
$ws.Cells.Item(73, 4).Value = 20235098
$ws.Cells.Item(74, 3).Value = 53917
$ws.Cells.Item(74, 4).Value = 78465154
$ws.Cells.Item(75, 3).Value = 151315
$ws.Cells.Item(75, 4).Value = 222909101
$ws.Cells.Item(76, 3).Value = 65419
$ws.Cells.Item(76, 4).Value = 97475418
$ws.Cells.Item(77, 3).Value = 20951
$ws.Cells.Item(77, 4).Value = 31306822
$ws.Cells.Item(78, 3).Value = 5013
$ws.Cells.Item(78, 4).Value = 7487903
$ws.Cells.Item(80, 3).Value = 23
$ws.Cells.Item(80, 4).Value = 33405
$ws.Cells.Item(85, 3).Value = 53146
$ws.Cells.Item(85, 4).Value = 72217207
$ws.Cells.Item(86, 3).Value = 4753
$ws.Cells.Item(86, 4).Value = 6888431
$ws.Cells.Item(87, 3).Value = 11847
$ws.Cells.Item(87, 4).Value = 17400043
$ws.Cells.Item(88, 3).Value = 3960
$ws.Cells.Item(88, 4).Value = 5900058
$ws.Cells.Item(89, 3).Value = 1368
$ws.Cells.Item(89, 4).Value = 2043789
$ws.Cells.Item(90, 3).Value = 294
$ws.Cells.Item(90, 4).Value = 438512
$ws.Cells.Item(93, 3).Value = 5549
$ws.Cells.Item(93, 4).Value = 7458190
$ws.Cells.Item(94, 3).Value = 1645
$ws.Cells.Item(94, 4).Value = 2369802
$ws.Cells.Item(95, 3).Value = 5310
$ws.Cells.Item(95, 4).Value = 7822137
$ws.Cells.Item(96, 3).Value = 1975
$ws.Cells.Item(96, 4).Value = 2940826
$ws.Cells.Item(101, 3).Value = 3661
$ws.Cells.Item(101, 4).Value = 4844503
$ws.Cells.Item(102, 3).Value = 702
$ws.Cells.Item(102, 4).Value = 1153875
$ws.Cells.Item(103, 3).Value = 423
$ws.Cells.Item(103, 4).Value = 711527
$ws.Cells.Item(104, 3).Value = 158
$ws.Cells.Item(104, 4).Value = 264520
$ws.Cells.Item(107, 3).Value = 11038
$ws.Cells.Item(107, 4).Value = 16012871
$ws.Cells.Item(108, 3).Value = 29684
$ws.Cells.Item(108, 4).Value = 43595433
$ws.Cells.Item(109, 3).Value = 9941
$ws.Cells.Item(109, 4).Value = 14780705
$ws.Cells.Item(110, 3).Value = 2743
$ws.Cells.Item(110, 4).Value = 4089580
$ws.Cells.Item(111, 3).Value = 503
$ws.Cells.Item(111, 4).Value = 749546
$ws.Cells.Item(112, 3).Value = 52
$ws.Cells.Item(112, 4).Value = 78000
$ws.Cells.Item(114, 3).Value = 9969
$ws.Cells.Item(114, 4).Value = 13163001
$ws.Cells.Item(115, 3).Value = 31119
$ws.Cells.Item(115, 4).Value = 44867671
$ws.Cells.Item(116, 3).Value = 67260
$ws.Cells.Item(116, 4).Value = 98417330
$ws.Cells.Item(117, 3).Value = 21690
$ws.Cells.Item(117, 4).Value = 32231988
$ws.Cells.Item(118, 3).Value = 6155
$ws.Cells.Item(118, 4).Value = 9169521
$ws.Cells.Item(119, 3).Value = 1155
$ws.Cells.Item(119, 4).Value = 1726100
$ws.Cells.Item(123, 3).Value = 5
$ws.Cells.Item(123, 4).Value = 7500
$ws.Cells.Item(124, 3).Value = 26272
$ws.Cells.Item(124, 4).Value = 35061310
$ws.Cells.Item(125, 3).Value = 36894
$ws.Cells.Item(125, 4).Value = 53237196
$ws.Cells.Item(126, 3).Value = 78314
$ws.Cells.Item(126, 4).Value = 114503319
$ws.Cells.Item(127, 3).Value = 24233
$ws.Cells.Item(127, 4).Value = 35965368
$ws.Cells.Item(128, 3).Value = 6519
$ws.Cells.Item(128, 4).Value = 9688358
$ws.Cells.Item(129, 3).Value = 1278
$ws.Cells.Item(129, 4).Value = 1900311
$ws.Cells.Item(130, 3).Value = 63
$ws.Cells.Item(130, 4).Value = 92728
$ws.Cells.Item(133, 3).Value = 32403
$ws.Cells.Item(133, 4).Value = 43003279
$ws.Cells.Item(134, 3).Value = 13562
$ws.Cells.Item(134, 4).Value = 19631470
$ws.Cells.Item(135, 3).Value = 32868
$ws.Cells.Item(135, 4).Value = 48266896
$ws.Cells.Item(136, 3).Value = 11649
$ws.Cells.Item(136, 4).Value = 17307587
$ws.Cells.Item(137, 3).Value = 3021
$ws.Cells.Item(137, 4).Value = 4502741
$ws.Cells.Item(138, 3).Value = 515
$ws.Cells.Item(138, 4).Value = 766490
$ws.Cells.Item(141, 3).Value = 10999
$ws.Cells.Item(141, 4).Value = 14657052
$ws.Cells.Item(142, 3).Value = 36002
$ws.Cells.Item(142, 4).Value = 51994212
$ws.Cells.Item(143, 3).Value = 83056
$ws.Cells.Item(143, 4).Value = 121678876
$ws.Cells.Item(144, 3).Value = 24793
$ws.Cells.Item(144, 4).Value = 36833619
$ws.Cells.Item(145, 3).Value = 6512
$ws.Cells.Item(145, 4).Value = 9717567
$ws.Cells.Item(146, 3).Value = 1477
$ws.Cells.Item(146, 4).Value = 2197730
$ws.Cells.Item(149, 3).Value = 29774
$ws.Cells.Item(149, 4).Value = 40137526
